$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row whose old B..G values it should receive
$rowMap = @{
  149 = 150
  150 = 149
  161 = 162
  162 = 163
  163 = 161
  264 = 265
  265 = 264
  316 = 318
  317 = 316
  318 = 317
  346 = 347
  347 = 346
  350 = 352
  351 = 350
  352 = 351
  355 = 356
  356 = 355
  375 = 376
  376 = 375
  379 = 380
  380 = 379
  389 = 390
  390 = 389
  400 = 401
  401 = 400
  419 = 420
  420 = 419
  431 = 432
  432 = 431
  536 = 537
  537 = 536
  583 = 584
  584 = 583
  586 = 587
  587 = 586
  590 = 591
  591 = 590
  593 = 594
  594 = 593
  687 = 688
  688 = 687
  709 = 710
  710 = 709
  720 = 721
  721 = 720
  889 = 890
  890 = 889
}

# Step 1: snapshot old values (columns B=2 .. G=7) for every affected row
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
  $rowVals = @{}
  foreach ($col in 2..7) {
    $rowVals[$col] = $ws.Cells.Item($r, $col).Value2
  }
  $snapshot[$r] = $rowVals
}

# Step 2: write new values using the snapshot (so overlapping writes never clobber unread data)
foreach ($r in $rowMap.Keys) {
  $src = $rowMap[$r]
  $srcVals = $snapshot[$src]
  foreach ($col in 2..7) {
    $ws.Cells.Item($r, $col).Value = $srcVals[$col]
  }
}
